$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Delete the two trailing duplicate rows (old "hasCopyrightResource" /
# "hasLicenseResource" rows 10-11) -- the replacement content below re-uses
# rows 2-9 for the full new "project metadata" property set, so the old
# rows 10-11 are no longer needed and the sheet shrinks by two rows.
# ---------------------------------------------------------------------------
$ws.Rows("10:11").Delete() | Out-Null

# ---------------------------------------------------------------------------
# Clear out the old metadata property rows (2-9) before writing the new
# project-metadata property rows in their place.
# ---------------------------------------------------------------------------
$ws.Range("A2:O9").ClearContents() | Out-Null

# Row 2: hasFileDescription (replaces hasCopyright)
$ws.Range("A2").Value = 'hasFileDescription'
$ws.Range("B2").Value = 'File Description'
$ws.Range("C2").Value = 'Dateibeschreibung'
$ws.Range("D2").Value = 'Description du fichier'
$ws.Range("E2").Value = 'Descrizione del file'
$ws.Range("L2").Value = 'hasValue, schema:description, crm:P190_has_symbolic_content'
$ws.Range("M2").Value = 'TextValue'
$ws.Range("N2").Value = 'Richtext'

# Row 3: hasFileName (replaces hasDescription)
$ws.Range("A3").Value = 'hasFileName'
$ws.Range("B3").Value = 'File Name'
$ws.Range("C3").Value = 'Dateiname'
$ws.Range("D3").Value = 'Nom de fichier'
$ws.Range("E3").Value = 'Nome del file'
$ws.Range("G3").Value = 'Name of the original file'
$ws.Range("H3").Value = 'Name der Originaldatei'
$ws.Range("I3").Value = 'Nom du fichier original'
$ws.Range("J3").Value = 'Nome del file originale'
$ws.Range("L3").Value = 'hasValue, schema:name'
$ws.Range("M3").Value = 'TextValue'
$ws.Range("N3").Value = 'SimpleText'

# Row 4: hasFileSize (was already hasFileSize, now gains comment text)
$ws.Range("A4").Value = 'hasFileSize'
$ws.Range("B4").Value = 'File Size (Mb)'
$ws.Range("C4").Value = 'Dateigrösse (Mb)'
$ws.Range("D4").Value = 'Taille du fichier (Mb)'
$ws.Range("E4").Value = 'Dimensioni del file (Mb)'
$ws.Range("G4").Value = 'Size of the file in Mb'
$ws.Range("H4").Value = 'Dateigröße in MB'
$ws.Range("I4").Value = 'Taille du fichier en Mo'
$ws.Range("J4").Value = 'Dimensione del file in Mb'
$ws.Range("L4").Value = 'hasValue, schema:size'
$ws.Range("M4").Value = 'DecimalValue'
$ws.Range("N4").Value = 'SimpleText'

# Row 5: hasID (was already hasID, now gains comment text)
$ws.Range("A5").Value = 'hasID'
$ws.Range("B5").Value = 'ID'
$ws.Range("C5").Value = 'ID'
$ws.Range("D5").Value = 'ID'
$ws.Range("E5").Value = 'ID'
$ws.Range("G5").Value = 'Unique identifier'
$ws.Range("H5").Value = 'Eindeutige Kennung'
$ws.Range("I5").Value = 'Identifiant unique'
$ws.Range("J5").Value = 'Identificatore univoco'
$ws.Range("L5").Value = 'hasValue, schema:identifier, crm:P1_is_identified_by'
$ws.Range("M5").Value = 'TextValue'
$ws.Range("N5").Value = 'SimpleText'

# Row 6: hasTimeStamp (replaces hasLicenseList)
$ws.Range("A6").Value = 'hasTimeStamp'
$ws.Range("B6").Value = 'Time Stamp'
$ws.Range("C6").Value = 'Zeitstempel'
$ws.Range("D6").Value = 'Horodatage'
$ws.Range("E6").Value = 'Timestamp'
$ws.Range("G6").Value = 'Time stamp'
$ws.Range("H6").Value = 'Zeitstempel'
$ws.Range("I6").Value = 'Horodatage'
$ws.Range("J6").Value = 'Data e ora'
$ws.Range("L6").Value = 'hasValue, schema:dateCreated'
$ws.Range("M6").Value = 'TimeValue'
$ws.Range("N6").Value = 'TimeStamp'

# Row 7: hasAuthorshipResource (replaces hasTimeStamp)
$ws.Range("A7").Value = 'hasAuthorshipResource'
$ws.Range("B7").Value = 'Author of the resource'
$ws.Range("C7").Value = 'Autor der Resource'
$ws.Range("D7").Value = 'Auteur·rice de la ressource'
$ws.Range("E7").Value = 'Autore della risorsa'
$ws.Range("G7").Value = 'Author of the resource'
$ws.Range("H7").Value = 'Autor der Resource'
$ws.Range("I7").Value = 'Auteur·rice de la ressource'
$ws.Range("J7").Value = 'Autore della risorsa'
$ws.Range("L7").Value = 'hasValue, foaf:person, schema:author, crm:E21_Person'
$ws.Range("M7").Value = 'TextValue'
$ws.Range("N7").Value = 'SimpleText'

# Row 8: hasCopyrightResource (moved up from old row 10; comments added)
$ws.Range("A8").Value = 'hasCopyrightResource'
$ws.Range("B8").Value = 'Copyright of the resource'
$ws.Range("C8").Value = 'Urheberrecht der Resource'
$ws.Range("D8").Value = 'Droits d''auteur de la ressource'
$ws.Range("E8").Value = 'Copyright della risorsa'
$ws.Range("G8").Value = 'Copyright of the resource'
$ws.Range("H8").Value = 'Urheberrecht der Resource'
$ws.Range("I8").Value = 'Droits d''auteur de la ressource'
$ws.Range("J8").Value = 'Copyright della risorsa'
$ws.Range("L8").Value = 'hasValue, schema:copyrightHolder, crm:P105_right_held_by'
$ws.Range("M8").Value = 'TextValue'
$ws.Range("N8").Value = 'SimpleText'

# Row 9: hasLicenseResource (moved up from old row 11; comments added)
$ws.Range("A9").Value = 'hasLicenseResource'
$ws.Range("B9").Value = 'License of the resource'
$ws.Range("C9").Value = 'Lizenz der Resource'
$ws.Range("D9").Value = 'Licence de la ressource'
$ws.Range("E9").Value = 'Licenza della risorsa'
$ws.Range("G9").Value = 'License of the resource'
$ws.Range("H9").Value = 'Lizenz der Resource'
$ws.Range("I9").Value = 'Licence de la ressource'
$ws.Range("J9").Value = 'Licenza della risorsa'
$ws.Range("L9").Value = 'hasValue, schema:license'
$ws.Range("M9").Value = 'ListValue'
$ws.Range("N9").Value = 'List'
$ws.Range("O9").Value = 'hlist: License'

# Re-select A4, matching the author's final cursor position.
$ws.Range("A4").Select() | Out-Null
